# Auto-generated from the commit diff: updates per-leve profit
# figures (currentAveragePrice*/LevePrice*/LeveProfit* columns)
# across several worksheets, as produced by the scheduled price-data
# refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 3415.5
$ws.Range("I12").Value = 750
$ws.Range("J12").Value = 4748.25
$ws.Range("K12").Value = 750
$ws.Range("L12").Value = 4748.25
$ws.Range("M12").Value = -580
$ws.Range("N12").Value = -5088.25

$ws.Range("H17").Value = 215923.92
$ws.Range("J17").Value = 215923.92
$ws.Range("L17").Value = 647771.76
$ws.Range("N17").Value = -648107.76

$ws.Range("H40").Value = 2163.5454
$ws.Range("I40").Value = 1700
$ws.Range("J40").Value = 2266.5557
$ws.Range("K40").Value = 1700
$ws.Range("L40").Value = 2266.5557
$ws.Range("M40").Value = -1525
$ws.Range("N40").Value = -2616.5557

$ws.Range("H112").Value = 1216.6666
$ws.Range("I112").Value = 610
$ws.Range("J112").Value = 1520
$ws.Range("K112").Value = 1830
$ws.Range("L112").Value = 4560
$ws.Range("M112").Value = -722
$ws.Range("N112").Value = -6776

$ws.Range("H132").Value = 5053.3335
$ws.Range("I132").Value = 1651.8518
$ws.Range("K132").Value = 4955.555399999999
$ws.Range("M132").Value = -2425.555399999999

$ws.Range("H137").Value = 1030553.8
$ws.Range("I137").Value = 1420.9354
$ws.Range("J137").Value = 2417646
$ws.Range("K137").Value = 4262.8062
$ws.Range("L137").Value = 7252938
$ws.Range("M137").Value = -1712.8062
$ws.Range("N137").Value = -7258038

$ws.Range("H138").Value = 2641.4243
$ws.Range("I138").Value = 1699.2559
$ws.Range("J138").Value = 4402.8696
$ws.Range("K138").Value = 5097.7677
$ws.Range("L138").Value = 13208.6088
$ws.Range("M138").Value = 42.23229999999967
$ws.Range("N138").Value = -23488.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1522.8387
$ws.Range("I2").Value = 1478.421
$ws.Range("J2").Value = 1593.1666
$ws.Range("K2").Value = 1478.421
$ws.Range("L2").Value = 1593.1666
$ws.Range("M2").Value = -1365.421
$ws.Range("N2").Value = -1819.1666

$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622

$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112

$ws.Range("H74").Value = 15552.177
$ws.Range("I74").Value = 18263.508
$ws.Range("J74").Value = 1502.5454
$ws.Range("K74").Value = 18263.508
$ws.Range("L74").Value = 1502.5454
$ws.Range("M74").Value = -17389.508
$ws.Range("N74").Value = -3250.5454

$ws.Range("H77").Value = 15552.177
$ws.Range("I77").Value = 18263.508
$ws.Range("J77").Value = 1502.5454
$ws.Range("K77").Value = 91317.54000000001
$ws.Range("L77").Value = 7512.727
$ws.Range("M77").Value = -86949.54000000001
$ws.Range("N77").Value = -16248.727

$ws.Range("H81").Value = 24000
$ws.Range("J81").Value = 24000
$ws.Range("L81").Value = 24000
$ws.Range("N81").Value = -25996

$ws.Range("H84").Value = 24000
$ws.Range("J84").Value = 24000
$ws.Range("L84").Value = 72000
$ws.Range("N84").Value = -81984

$ws.Range("H102").Value = 1870
$ws.Range("I102").Value = 1870
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1870
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -248
$ws.Range("N102").ClearContents()

$ws.Range("H113").Value = 30000
$ws.Range("J113").Value = 30000
$ws.Range("L113").Value = 30000
$ws.Range("N113").Value = -38678

$ws.Range("H116").Value = 1522.8387
$ws.Range("I116").Value = 1478.421
$ws.Range("J116").Value = 1593.1666
$ws.Range("K116").Value = 1478.421
$ws.Range("L116").Value = 1593.1666
$ws.Range("M116").Value = 815.579
$ws.Range("N116").Value = -6181.1666

$ws.Range("H132").Value = 1903.9286
$ws.Range("I132").Value = 1229.1428
$ws.Range("K132").Value = 3687.4284
$ws.Range("M132").Value = -1157.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1522.8387
$ws.Range("I3").Value = 1478.421
$ws.Range("J3").Value = 1593.1666
$ws.Range("K3").Value = 1478.421
$ws.Range("L3").Value = 1593.1666
$ws.Range("M3").Value = -1364.421
$ws.Range("N3").Value = -1821.1666

$ws.Range("H22").Value = 840.2308
$ws.Range("I22").Value = 853.5833
$ws.Range("J22").Value = 680
$ws.Range("K22").Value = 853.5833
$ws.Range("L22").Value = 680
$ws.Range("M22").Value = -680.5833
$ws.Range("N22").Value = -1026

$ws.Range("H99").Value = 1277.2778
$ws.Range("I99").Value = 944.625
$ws.Range("J99").Value = 1543.4
$ws.Range("K99").Value = 944.625
$ws.Range("L99").Value = 1543.4
$ws.Range("M99").Value = 553.375
$ws.Range("N99").Value = -4539.4

$ws.Range("H105").Value = 2504.4
$ws.Range("J105").Value = 2500
$ws.Range("L105").Value = 2500
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1272
$ws.Range("I47").Value = 566.6667
$ws.Range("J47").Value = 1801
$ws.Range("K47").Value = 1700.0001
$ws.Range("L47").Value = 5403
$ws.Range("M47").Value = -1269.0001
$ws.Range("N47").Value = -6265

$ws.Range("H102").Value = 5874.25
$ws.Range("J102").Value = 7665.6665
$ws.Range("L102").Value = 22996.9995
$ws.Range("N102").Value = -27864.9995

$ws.Range("H131").Value = 982.87
$ws.Range("J131").Value = 982.87
$ws.Range("L131").Value = 2948.61
$ws.Range("N131").Value = -13028.61

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 868.46155
$ws.Range("I16").Value = 887.7778
$ws.Range("K16").Value = 887.7778
$ws.Range("M16").Value = -717.7778

$ws.Range("H22").Value = 46140.273
$ws.Range("I22").Value = 250420
$ws.Range("J22").Value = 744.7778
$ws.Range("K22").Value = 250420
$ws.Range("L22").Value = 744.7778
$ws.Range("M22").Value = -250125
$ws.Range("N22").Value = -1334.7778

$ws.Range("H27").Value = 46140.273
$ws.Range("I27").Value = 250420
$ws.Range("J27").Value = 744.7778
$ws.Range("K27").Value = 250420
$ws.Range("L27").Value = 744.7778
$ws.Range("M27").Value = -250313
$ws.Range("N27").Value = -958.7778

$ws.Range("H68").Value = 17699.54
$ws.Range("I68").Value = 19066.166
$ws.Range("J68").Value = 1300
$ws.Range("K68").Value = 19066.166
$ws.Range("L68").Value = 1300
$ws.Range("M68").Value = -18317.166
$ws.Range("N68").Value = -2798

$ws.Range("H71").Value = 17699.54
$ws.Range("I71").Value = 19066.166
$ws.Range("J71").Value = 1300
$ws.Range("K71").Value = 95330.83
$ws.Range("L71").Value = 6500
$ws.Range("M71").Value = -91586.83
$ws.Range("N71").Value = -13988

$ws.Range("H82").Value = 975
$ws.Range("I82").Value = 900
$ws.Range("J82").Value = 981.8182
$ws.Range("K82").Value = 900
$ws.Range("L82").Value = 981.8182
$ws.Range("M82").Value = -539
$ws.Range("N82").Value = -1703.8182

$ws.Range("H85").Value = 975
$ws.Range("I85").Value = 900
$ws.Range("J85").Value = 981.8182
$ws.Range("K85").Value = 900
$ws.Range("L85").Value = 981.8182
$ws.Range("M85").Value = 348
$ws.Range("N85").Value = -3477.8182

$ws.Range("H100").Value = 1776.7059
$ws.Range("I100").Value = 1443.4286
$ws.Range("K100").Value = 1443.4286
$ws.Range("M100").Value = -902.4286

$ws.Range("H132").Value = 6193.154
$ws.Range("I132").Value = 7574.7144
$ws.Range("J132").Value = 4581.3335
$ws.Range("K132").Value = 22724.1432
$ws.Range("L132").Value = 13744.0005
$ws.Range("M132").Value = -20194.1432
$ws.Range("N132").Value = -18804.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5779.7896
$ws.Range("I107").Value = 12888.125
$ws.Range("J107").Value = 610.0909
$ws.Range("K107").Value = 38664.375
$ws.Range("L107").Value = 1830.2727
$ws.Range("M107").Value = -36744.375
$ws.Range("N107").Value = -5670.2727

$ws.Range("H136").Value = 3244.575
$ws.Range("I136").Value = 496.45456
$ws.Range("J136").Value = 16200
$ws.Range("K136").Value = 1489.36368
$ws.Range("L136").Value = 48600
$ws.Range("M136").Value = 1060.63632
$ws.Range("N136").Value = -53700

$excel.Calculate()

